$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps storing values as text, matching the
# original inline-string cell type (many look numeric, e.g. "43.058.10").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.058.10'
$ws.Range("E2").Value = '  +0.88%  '
$ws.Range("D3").Value = '2.350.45'
$ws.Range("E3").Value = '  +4.75%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '311.32'
$ws.Range("E5").Value = '  +4.99%  '
$ws.Range("D6").Value = '108.35'
$ws.Range("E6").Value = '  -4.45%  '
$ws.Range("D7").Value = '0.635'
$ws.Range("E7").Value = '  +0.86%  '
$ws.Range("E8").Value = '  -0.46%  '
$ws.Range("D9").Value = '0.624'
$ws.Range("D10").Value = '43.53'
$ws.Range("E10").Value = '  -4.81%  '
$ws.Range("D11").Value = '0.0938'
$ws.Range("E11").Value = '  +1.09%  '
$ws.Range("D12").Value = '8.95'
$ws.Range("E12").Value = '  -1.15%  '
$ws.Range("D13").Value = '1.08'
$ws.Range("E13").Value = '  +18.67%  '
$ws.Range("E14").Value = '  +0.44%  '
$ws.Range("D15").Value = '16.37'
$ws.Range("E15").Value = '  +6.94%  '
$ws.Range("D16").Value = '2.694.79'
$ws.Range("E16").Value = '  +4.32%  '
$ws.Range("D17").Value = '2.408.72'
$ws.Range("E17").Value = '  +6.74%  '
$ws.Range("D18").Value = '42.991.44'
$ws.Range("E18").Value = '  +0.52%  '
$ws.Range("E19").Value = '  +0.16%  '
$ws.Range("D20").Value = '7.24'
$ws.Range("E20").Value = '  -4.00%  '
$ws.Range("D21").Value = '75.67'
$ws.Range("E21").Value = '  +2.98%  '
$ws.Range("D22").Value = '3.47'
$ws.Range("E22").Value = '  -2.59%  '
$ws.Range("E23").Value = '  +7.64%  '
$ws.Range("D24").Value = '249.71'
$ws.Range("E24").Value = '  +7.73%  '
$ws.Range("D25").Value = '8.95'
$ws.Range("E25").Value = '  -5.92%  '
$ws.Range("D26").Value = '11.91'
$ws.Range("E26").Value = '  -1.09%  '
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("D28").Value = '2.24'
$ws.Range("E28").Value = '  +0.27%  '
$ws.Range("D29").Value = '38.73'
$ws.Range("E29").Value = '  -3.24%  '
$ws.Range("D30").Value = '22.52'
$ws.Range("E30").Value = '  +5.73%  '
$ws.Range("D31").Value = '174.13'
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("D32").Value = '3.17'
$ws.Range("E32").Value = '  -2.90%  '
$ws.Range("D33").Value = '0.0909'
$ws.Range("E33").Value = '  +1.07%  '
$ws.Range("E34").Value = '  +0.98%  '
$ws.Range("D35").Value = '4.97'
$ws.Range("E35").Value = '  -1.33%  '
$ws.Range("D36").Value = '0.131'
$ws.Range("E36").Value = '  +1.93%  '
$ws.Range("D37").Value = '0.0378'
$ws.Range("E37").Value = '  +1.65%  '
$ws.Range("D38").Value = '4.13'
$ws.Range("E38").Value = '  -4.57%  '
$ws.Range("E39").Value = '  -1.47%  '
$ws.Range("D40").Value = '2.78'
$ws.Range("E40").Value = '  +8.65%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = '1.50'
$ws.Range("E41").Value = '  +12.43%  '
$ws.Range("B42").Value = 'MultiversX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D42").Value = '72.02'
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("E43").Value = '  -3.00%  '
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("D45").Value = '12.46'
$ws.Range("E45").Value = '  -6.52%  '
$ws.Range("D46").Value = '5.70'
$ws.Range("E46").Value = '  +2.10%  '
$ws.Range("D47").Value = '9.25'
$ws.Range("E47").Value = '  +6.50%  '
$ws.Range("D48").Value = '110.24'
$ws.Range("E48").Value = '  +3.31%  '
$ws.Range("D49").Value = '1.30'
$ws.Range("E49").Value = '  -1.44%  '
$ws.Range("D50").Value = '0.0999'
$ws.Range("E50").Value = '  +0.97%  '
$ws.Range("D51").Value = '70.73'
$ws.Range("E51").Value = '  +3.60%  '
